# Update the PRISMA flowchart numbers for the "papers up to 2023" refresh.
# The flowchart lives inside a single group shape ("Group 50") on slide 1;
# each box is a TextBox shape inside that group (Shape.GroupItems).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)

function Set-BoxText($Group, $ItemIndex, $ParaIndex, $NewText) {
    $shape = $Group.GroupItems.Item($ItemIndex)
    $tr = $shape.TextFrame.TextRange
    $para = $tr.Paragraphs($ParaIndex, 1)
    $para.Text = $NewText
}

# TextBox 7 (id=8): "Records from databases (...)" / n = 1829 -> n = 3070
Set-BoxText $grp 4 3 "n = 3070"

# TextBox 8 (id=9): "Records after duplicates removed" / n = 1591 -> n = 2675
Set-BoxText $grp 5 3 "n = 2675"

# TextBox 9 (id=10): "Records screened" / n = 1591 -> n = 2675
Set-BoxText $grp 6 3 "n = 2675"

# TextBox 15 (id=16): "Records excluded based on title and abstract (n = 1331)" -> (n = 2316)
Set-BoxText $grp 9 1 "Records excluded based on title and abstract (n = 2316)"

# TextBox 19 (id=20): "Full text articles assessed for eligibility" / n = 260 -> n = 359
Set-BoxText $grp 11 3 "n = 359"

# TextBox 25 (id=26): list of exclusion reasons with counts
Set-BoxText $grp 12 1 "Full text articles excluded for following reasons (n = 228):"
Set-BoxText $grp 12 2 "Not an implemented AR application (n = 81)"
Set-BoxText $grp 12 3 "Not for education (n = 60)"
Set-BoxText $grp 12 4 "Not interactive, collaborative or multiuser (n = 47)"
Set-BoxText $grp 12 5 "Outside target audience (n = 36)"
Set-BoxText $grp 12 6 "Not peer reviewed (n = 4)"

# TextBox 27 (id=28): "Studies included in the literature review" / n = 100 -> n = 131
Set-BoxText $grp 13 3 "n = 131"
